$wb = $excel.ActiveWorkbook

# The two data sheets ("展览" and "全部类型") carry identical content and
# both need the same cell updates.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 24
    $ws.Range("F3").Value = 1840
    $ws.Range("F5").Value = 801
    $ws.Range("F10").Value = 241
    $ws.Range("F15").Value = 4316
    $ws.Range("F18").Value = 470
    $ws.Range("F21").Value = 1647
    $ws.Range("I21").Value = "//i2.hdslb.com/bfs/openplatform/202404/S1nqZf721712025221477.jpeg"
    $ws.Range("F22").Value = 364
    $ws.Range("F24").Value = 8
    $ws.Range("F26").Value = 2016
    $ws.Range("F29").Value = 4
    $ws.Range("F30").Value = 140
    $ws.Range("F31").Value = 66
    $ws.Range("F32").Value = 206
}
